$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.094.57"
$ws.Range("E2").Value = "  +2.05%  "

$ws.Range("D3").Value = "1.909.36"
$ws.Range("E3").Value = "  +1.80%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -1.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.92"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.96%  "

$ws.Range("E7").Value = "  +0.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3823"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07356"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9325"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.79"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07817"
$ws.Range("D12").ClearFormats()

$ws.Range("D13").Value = "1.904.41"
$ws.Range("E13").Value = "  +1.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.507"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.629"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.34"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008825"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.90%  "

$ws.Range("D20").Value = "28.101.83"
$ws.Range("E20").Value = "  +1.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.83"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.150"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("D23").Value = "2.131.45"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.88"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.94"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.913"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.58"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("E28").Value = "  +4.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.40"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.959"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08910"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.343"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.251"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.88%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7674"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.36%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.682"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.615"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02048"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.101"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05300"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5484"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.976"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.023"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1523"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.460"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.73"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4823"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "107.54"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.656"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.31"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06104"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.14%  "
